$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): columns reordered + new audit columns appended ----
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "method_seq"
$ws.Range("C1").Value = "lang_code"
$ws.Range("D1").Value = "is_active"
$ws.Range("E1").Value = "cr_by"
$ws.Range("F1").Value = "cr_dtimes"
$ws.Range("G1").Value = "upd_by"
$ws.Range("H1").Value = "upd_dtimes"
$ws.Range("I1").Value = "is_deleted"
$ws.Range("J1").Value = "del_dtimes"

# ---- Apply the date/time display format to column F (cr_dtimes) ----
$ws.Range("F2:F6").NumberFormat = "mm:ss.0"

# ---- Data rows ----
$codes = "PWD", "OTP", "FINGERPRINT", "IRIS", "FACE"
for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $codes[$i]
    $ws.Cells.Item($r, 2).Value = $i + 1
    $ws.Cells.Item($r, 3).Value = "fra"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = 45079.576795277775
    $ws.Cells.Item($r, 7).Value = "NULL"
    $ws.Cells.Item($r, 8).Value = "NULL"
    $ws.Cells.Item($r, 9).Value = $false
    $ws.Cells.Item($r, 10).Value = "NULL"
}

$ws.Range("F13").Select() | Out-Null
